$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("M5").Value = 25413.36000000007
$ws.Range("P5").Value = 0.04025

# Row 38
$ws.Range("C38").Value = 415
$ws.Range("H38").Value = 15
$ws.Range("I38").Value = 59.4
$ws.Range("J38").Value = 29.4
$ws.Range("M38").Value = 24592.80000000015
$ws.Range("N38").Value = 95
$ws.Range("O38").Value = 0.0019
$ws.Range("P38").Value = 0.04113

# Row 39
$ws.Range("C39").Value = 433
$ws.Range("H39").Value = 33
$ws.Range("I39").Value = 58.68
$ws.Range("J39").Value = 28.68
$ws.Range("M39").Value = 24586.20000000011
$ws.Range("N39").Value = 90
$ws.Range("O39").Value = 0.0018
$ws.Range("P39").Value = 0.03871

# Row 44
$ws.Range("O44").Value = 0.0019
$ws.Range("P44").Value = 0.0413

# Row 46
$ws.Range("O46").Value = 0.0024
$ws.Range("P46").Value = 0.05217

# Row 47
$ws.Range("M47").Value = 24413.28000000015
$ws.Range("P47").Value = 0.03888

# Row 48
$ws.Range("O48").Value = 0.0025
$ws.Range("P48").Value = 0.05447

# Row 49
$ws.Range("O49").Value = 0.0019
$ws.Range("P49").Value = 0.04121

# Row 50
$ws.Range("C50").Value = 439
$ws.Range("H50").Value = 39
$ws.Range("I50").Value = 58.44
$ws.Range("J50").Value = 28.44
$ws.Range("M50").Value = 24240.83999999985
$ws.Range("O50").Value = 0.0024
$ws.Range("P50").Value = 0.05206

# Row 51
$ws.Range("C51").Value = 438
$ws.Range("H51").Value = 38
$ws.Range("I51").Value = 58.48
$ws.Range("J51").Value = 28.48
$ws.Range("M51").Value = 24200.79999999984
$ws.Range("O51").Value = 0.0019
$ws.Range("P51").Value = 0.0413
